# Auto-generated Excel COM-interop script to update country case data
# and append two new country rows (San Cristobal y Nieves, Islas Virgenes Britanicas)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 25 de Marzo de 2020 a las 23:16'

$data = @(
  @('China',81218,0,73650,4287,1399,0,3281),
  @('Italia',74386,5210,9362,57521,3489,683,7503),
  @('Estados Unidos',64832,9976,393,63526,1411,133,913),
  @('España',49515,7457,5367,40501,3166,656,3647),
  @('Alemania',37323,4332,3547,33570,23,47,206),
  @('Iran',27017,2206,9625,15315,0,143,2077),
  @('Francia',25233,2929,3900,20002,2827,231,1331),
  @('Suiza',10897,1020,131,10613,141,31,153),
  @('Reino Unido',9529,1452,135,8929,163,43,465),
  @('Corea del Sur',9137,100,3730,5281,59,6,126),
  @('Paises Bajos',6412,852,3,6053,582,80,356),
  @('Austria',5588,305,9,5549,28,2,30),
  @('Belgica',4937,668,547,4212,474,56,178),
  @('Canada',3367,575,185,3152,1,4,30),
  @('Noruega',3066,200,6,3046,57,2,14),
  @('Portugal',2995,633,22,2930,61,10,43),
  @('Suecia',2526,227,16,2448,158,22,62),
  @('Brasil',2433,186,2,2374,18,11,57),
  @('Turquia',2433,561,26,2348,136,15,59),
  @('Australia',2431,114,118,2304,11,1,9),
  @('Israel',2369,439,58,2306,37,2,5),
  @('Malasia',1796,172,199,1577,64,4,20),
  @('Dinamarca',1724,133,1,1689,69,2,34),
  @('Chequia',1654,260,10,1638,19,3,6),
  @('Irlanda',1564,235,5,1550,39,2,9),
  @('Luxemburgo',1333,234,6,1319,3,0,8),
  @('Japon',1307,114,310,952,57,2,45),
  @('Ecuador',1173,91,3,1142,2,1,28),
  @('Chile',1142,220,22,1117,7,1,3),
  @('Pakistan',1063,91,21,1034,5,1,8),
  @('Polonia',1051,150,7,1030,3,4,14),
  @('Tailandia',934,107,70,860,11,0,4),
  @('Rumania',906,112,86,806,18,2,14),
  @('Arabia Saudita',900,133,29,869,0,1,2),
  @('Finlandia',880,88,10,867,22,2,3),
  @('Grecia',821,78,36,763,53,2,22),
  @('Indonesia',790,104,31,701,0,3,58),
  @('Islandia',737,89,56,679,11,0,2),
  @('Crucero',712,0,587,115,15,0,10),
  @('Sudafrica',709,155,12,697,2,0,0),
  @('Rusia',658,163,29,626,8,2,3),
  @('India',657,121,43,602,0,2,12),
  @('Filipinas',636,84,26,572,1,3,38),
  @('Singapur',631,73,160,469,17,0,2),
  @('Catar',537,11,41,496,6,0,0),
  @('Eslovenia',528,48,10,513,14,1,5),
  @('Peru',480,64,1,470,9,2,9),
  @('Colombia',470,92,8,458,0,1,4),
  @('Egipto',456,54,95,340,0,1,21),
  @('Panama',443,0,1,434,33,2,8),
  @('Croacia',442,60,22,419,6,0,1),
  @('Barein',419,27,177,238,2,1,4),
  @('Hong Kong',410,23,102,304,4,0,4),
  @('Mexico',405,38,4,396,1,1,5),
  @('Estonia',404,35,8,395,6,1,1),
  @('Republica Dominicana',392,80,3,379,0,4,10),
  @('Argentina',387,0,52,327,0,2,8),
  @('Serbia',384,81,15,365,21,1,4),
  @('Irak',346,30,103,214,0,2,29),
  @('Libano',333,15,20,307,4,2,6),
  @('Emiratos Arabes Unidos',333,85,52,279,2,0,2),
  @('Argelia',302,38,65,216,0,2,21),
  @('Lituania',274,65,1,269,1,2,4),
  @('Armenia',265,16,16,249,6,0,0),
  @('Bulgaria',242,24,4,235,8,0,3),
  @('Taiwan',235,19,29,204,0,0,2),
  @('Hungria',226,39,21,195,6,1,10),
  @('Marruecos',225,55,7,212,1,1,6),
  @('Letonia',221,24,1,220,0,0,0),
  @('Eslovaquia',216,12,7,209,2,0,0),
  @('San Marino',208,21,4,183,12,0,21),
  @('Nueva Zelanda',205,0,22,183,0,0,0),
  @('Costa Rica',201,24,2,197,4,0,2),
  @('Kuwait',195,4,43,152,6,0,0),
  @('Uruguay',189,0,0,189,3,0,0),
  @('Principado de Andorra',188,24,1,186,6,0,1),
  @('Republica de Macedonia',177,29,1,173,1,1,3),
  @('Bosnia y Herzegovina',176,8,2,171,1,0,3),
  @('Tunez',173,59,2,166,11,1,5),
  @('Jordania',172,18,1,171,0,0,0),
  @('Moldavia',149,24,2,146,20,0,1),
  @('Burkina Faso',146,32,10,132,0,0,4),
  @('Albania',146,23,17,124,3,0,5),
  @('Ucrania',145,43,1,139,0,2,5),
  @('Vietnam',141,7,17,124,3,0,0),
  @('Republica de Chipre',132,8,3,126,3,0,3),
  @('Islas Feroe',132,10,38,94,2,0,0),
  @('Malta',129,19,2,127,1,0,0),
  @('Reunion',111,17,1,110,0,0,0),
  @('Brunei',109,5,2,107,1,0,0),
  @('Sri Lanka',102,0,3,99,2,0,0),
  @('Senegal',99,13,9,90,0,0,0),
  @('Oman',99,15,17,82,0,0,0),
  @('Camboya',96,5,10,86,1,0,0),
  @('Ghana',93,40,0,89,0,2,4),
  @('Azerbaiyan',93,6,10,81,6,1,2),
  @('Venezuela',91,7,15,76,2,0,0),
  @('Bielorrusia',86,5,29,57,2,0,0),
  @('Afganistan',84,10,2,80,0,1,2),
  @('Kazajistan',81,9,0,81,0,0,0),
  @('Camerun',75,9,2,72,0,0,1),
  @('Georgia',75,5,10,65,1,0,0),
  @('Guadalupe',73,0,0,72,4,0,1),
  @('Costa de Marfil',73,0,3,70,0,0,0),
  @('Estado de Palestina',71,11,16,54,0,1,1),
  @('Martinica',66,9,0,65,7,0,1),
  @('Uzbekistan',60,10,0,60,4,0,0),
  @('Trinidad yTobago',60,3,0,60,0,0,0),
  @('Cuba',57,9,1,55,2,0,1),
  @('Montenegro',52,5,0,51,0,0,1),
  @('Liechtenstein',51,0,0,51,0,0,0),
  @('Consejo Danes para los Refugiados',48,3,0,46,0,0,2),
  @('Mauricio',48,6,0,46,1,0,2),
  @('Nigeria',46,2,2,43,0,0,1),
  @('Kirguistan',44,2,0,44,0,0,0),
  @('Ruanda',41,1,0,41,0,0,0),
  @('Puerto Rico',39,0,1,36,0,0,2),
  @('Banglades',39,0,7,27,0,1,5),
  @('Paraguay',37,10,0,34,1,1,3),
  @('Honduras',36,6,0,36,0,0,0),
  @('Mayotte',36,0,0,36,0,0,0),
  @('Bolivia',32,3,0,32,0,0,0),
  @('Guam',32,0,0,31,0,0,1),
  @('Monaco',31,8,1,30,0,0,0),
  @('Macao',30,1,10,20,0,0,0),
  @('Kenia',28,3,1,27,0,0,0),
  @('Guayana Francesa',28,5,6,22,0,0,0),
  @('Gibraltar',26,11,5,21,0,0,0),
  @('Polinesia Francesa',25,0,0,25,0,0,0),
  @('Jamaica',25,4,2,22,0,0,1),
  @('Guatemala',24,3,4,19,0,0,1),
  @('Isla de Man',23,0,0,23,0,0,0),
  @('Togo',23,3,1,22,0,0,0),
  @('Madagascar',19,2,0,19,0,0,0),
  @('Barbados',18,0,0,18,0,0,0),
  @('Islas Virgenes de los Estados Unidos',17,0,0,17,0,0,0),
  @('Aruba',17,0,1,16,0,0,0),
  @('Nueva Caledonia',14,4,0,14,0,0,0),
  @('Uganda',14,5,0,14,0,0,0),
  @('Maldivas',13,0,8,5,0,0,0),
  @('Zambia',12,9,0,12,0,0,0),
  @('Etiopia',12,0,0,12,0,0,0),
  @('Tanzania',12,0,0,12,0,0,0),
  @('Republica de Yibuti',11,8,0,11,0,0,0),
  @('San Martin (Parte Francesa)',11,3,0,11,0,0,0),
  @('Mongolia',10,0,0,10,0,0,0),
  @('El Salvador',9,4,0,9,0,0,0),
  @('Guinea Ecuatorial',9,0,0,9,0,0,0),
  @('Haiti',8,1,0,8,0,0,0),
  @('Surinam',8,1,0,8,0,0,0),
  @('Islas Caimanes',8,2,0,7,0,0,1),
  @('Seychelles',7,0,0,7,0,0,0),
  @('Bermudas',7,1,0,7,0,0,0),
  @('Dominica',7,0,0,7,0,0,0),
  @('Niger',7,4,0,6,0,1,1),
  @('Namibia',7,0,2,5,0,0,0),
  @('Benin',6,0,0,6,0,0,0),
  @('Curazao',6,0,0,5,0,0,1),
  @('Gabon',6,0,0,5,0,0,1),
  @('Groenlandia',6,1,2,4,0,0,0),
  @('Mozambique',5,2,0,5,0,0,0),
  @('Siria',5,4,0,5,0,0,0),
  @('Fiyi',5,1,0,5,0,0,0),
  @('Bahamas',5,0,1,4,0,0,0),
  @('Guyana',5,0,0,4,0,0,1),
  @('Eritrea',4,3,0,4,0,0,0),
  @('Santa Sede',4,0,0,4,0,0,0),
  @('Congo',4,0,0,4,0,0,0),
  @('Guinea',4,0,0,4,0,0,0),
  @('Suazilandia',4,0,0,4,0,0,0),
  @('Cabo Verde',4,1,0,3,0,0,1),
  @('Liberia',3,0,0,3,0,0,0),
  @('San Martin (Parte Holandesa)',3,1,0,3,0,0,0),
  @('Angola',3,0,0,3,0,0,0),
  @('Antigua y Barbuda',3,0,0,3,0,0,0),
  @('Birmania',3,0,0,3,0,0,0),
  @('San Bartolome',3,0,0,3,0,0,0),
  @('Republica del Chad',3,0,0,3,0,0,0),
  @('Republica de Africa Central',3,0,0,3,0,0,0),
  @('Santa Lucia',3,0,0,3,0,0,0),
  @('Laos',3,1,0,3,0,0,0),
  @('Gambia',3,0,0,2,0,0,1),
  @('Sudan',3,0,0,2,0,0,1),
  @('Nepal',3,1,1,2,0,0,0),
  @('Zimbabue',3,0,0,2,0,0,1),
  @('Guinea-Bisau',2,2,0,2,0,0,0),
  @('San Cristobal y Nieves',2,2,0,2,0,0,0),
  @('Butan',2,0,0,2,0,0,0),
  @('Nicaragua',2,0,0,2,0,0,0),
  @('Mauritania',2,0,0,2,0,0,0),
  @('Islas Virgenes Britanicas',2,2,0,2,0,0,0),
  @('Mali',2,2,0,2,0,0,0),
  @('San Vicente y las Granadinas',1,0,0,1,0,0,0),
  @('Granada',1,0,0,1,0,0,0),
  @('Islas Turcas y Caicos',1,0,0,1,0,0,0),
  @('Belice',1,0,0,1,0,0,0),
  @('Somalia',1,0,0,1,0,0,0),
  @('Papua Nueva Guinea',1,0,0,1,0,0,0),
  @('Montserrat',1,0,0,1,0,0,0),
  @('Libia',1,0,0,1,0,0,0),
  @('Timor Oriental',1,0,0,1,0,0,0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
}

Write-Output "Done updating $($data.Length) data rows"
